# Finish Goods Inward Template - add a "finishedGoodsTotalQty" column
# between "fgTotalQty" (I) and "finishedGoodsUnit" (old J), and bump the
# sample date forward by a month.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new column at J -------------------------------------------
# This shifts the old J:M ("finishedGoodsUnit", "scrapQty", "scrapUnit",
# "finishedBy") one position to the right (K:N) and keeps their widths and
# styles intact.
$ws.Columns("J").Insert()

# --- New header (row 1) ---------------------------------------------------
# The inserted column already carries the header row's wrap-text style,
# so only the value needs to be set.
$ws.Cells.Item(1, 10).Value = "finishedGoodsTotalQty"

# --- New sample data (row 2) ----------------------------------------------
$ws.Cells.Item(2, 10).Value = 100

# --- Update the sample date in A2 (2020-06-25 -> 2020-07-25) -------------
$ws.Cells.Item(2, 1).Value = 44037

# --- Column J should visually match column I's width (the column it was
# inserted next to) ---------------------------------------------------------
$ws.Columns("J").ColumnWidth = $ws.Columns("I").ColumnWidth

# --- Update the view: scroll so column E is the left-most visible column
# and select the new header cell J1 ----------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 5
[void]$ws.Range("J1").Select()

Write-Host "Inserted finishedGoodsTotalQty column and updated sample row."
